# Adds daily report rows for Montese from 2021-09-21 (row 386) through
# 2021-12-08 (row 464): "aggiornamento fino a 8/12".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(44460,0,4,122.1374045801527),
    @(44461,0,4,122.1374045801527),
    @(44462,0,4,122.1374045801527),
    @(44463,0,4,122.1374045801527),
    @(44464,1,3,91.6030534351145),
    @(44465,0,1,30.53435114503817),
    @(44466,0,1,30.53435114503817),
    @(44467,0,1,30.53435114503817),
    @(44468,0,1,30.53435114503817),
    @(44469,0,1,30.53435114503817),
    @(44470,0,1,30.53435114503817),
    @(44471,0,0,0),
    @(44472,0,0,0),
    @(44473,0,0,0),
    @(44474,0,0,0),
    @(44475,0,0,0),
    @(44476,0,0,0),
    @(44477,0,0,0),
    @(44478,0,0,0),
    @(44479,0,0,0),
    @(44480,0,0,0),
    @(44481,0,0,0),
    @(44482,0,0,0),
    @(44483,0,0,0),
    @(44484,0,0,0),
    @(44485,0,0,0),
    @(44486,0,0,0),
    @(44487,0,0,0),
    @(44488,0,0,0),
    @(44489,0,0,0),
    @(44490,0,0,0),
    @(44491,0,0,0),
    @(44492,0,0,0),
    @(44493,0,0,0),
    @(44494,0,0,0),
    @(44495,0,0,0),
    @(44496,0,0,0),
    @(44497,0,0,0),
    @(44498,0,0,0),
    @(44499,0,0,0),
    @(44500,0,0,0),
    @(44501,0,0,0),
    @(44502,0,0,0),
    @(44503,0,0,0),
    @(44504,0,0,0),
    @(44505,0,0,0),
    @(44506,1,1,30.53435114503817),
    @(44507,0,1,30.53435114503817),
    @(44508,1,2,61.06870229007634),
    @(44509,0,2,61.06870229007634),
    @(44510,0,2,61.06870229007634),
    @(44511,0,2,61.06870229007634),
    @(44512,0,2,61.06870229007634),
    @(44513,0,1,30.53435114503817),
    @(44514,0,1,30.53435114503817),
    @(44515,0,0,0),
    @(44516,5,5,152.6717557251908),
    @(44517,1,6,183.206106870229),
    @(44518,0,6,183.206106870229),
    @(44519,0,6,183.206106870229),
    @(44520,0,6,183.206106870229),
    @(44521,0,6,183.206106870229),
    @(44522,1,7,213.7404580152672),
    @(44523,0,2,61.06870229007634),
    @(44524,5,6,183.206106870229),
    @(44525,1,7,213.7404580152672),
    @(44526,0,7,213.7404580152672),
    @(44527,0,7,213.7404580152672),
    @(44528,0,7,213.7404580152672),
    @(44529,0,6,183.206106870229),
    @(44530,0,6,183.206106870229),
    @(44531,0,1,30.53435114503817),
    @(44532,1,1,30.53435114503817),
    @(44533,0,1,30.53435114503817),
    @(44534,1,2,61.06870229007634),
    @(44535,3,5,152.6717557251908),
    @(44536,4,9,274.8091603053435),
    @(44537,3,12,366.412213740458),
    @(44538,0,12,366.412213740458)
)

$startRow = 386
$endRow = $startRow + $data.Count - 1

# Column A keeps the date-style (border + bold + centered + numFmt 165)
# used throughout the rest of the sheet, so copy that formatting down
# before writing the new values.
$ws.Range("A385").Copy() | Out-Null
$ws.Range("A$startRow`:A$endRow").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$colA  = New-Object 'object[,]' $data.Count,1
$colBD = New-Object 'object[,]' $data.Count,3

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $colA[$i,0]  = $row[0]
    $colBD[$i,0] = $row[1]
    $colBD[$i,1] = $row[2]
    $colBD[$i,2] = $row[3]
}

$ws.Range("A$startRow`:A$endRow").Value = $colA
$ws.Range("B$startRow`:D$endRow").Value = $colBD
